$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "harbor_type"
$ws.Range("D1").Value = "harbor_resource"

# Data: harbor_type (C) and harbor_resource (D) for rows 2-19
$data = @(
    @(2, 5),
    @(2, 5),
    @(1, 3),
    @(1, 3),
    @(1, 2),
    @(1, 2),
    @(2, 5),
    @(2, 5),
    @(1, 4),
    @(1, 4),
    @(2, 5),
    @(2, 5),
    @(2, 5),
    @(2, 5),
    @(1, 0),
    @(1, 0),
    @(1, 1),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}

# Update selection to match target
$ws.Range("H4").Select()
